# ---------------------------------------------------------------------------
# historico.xlsx — "Add files via upload"
#
# The leaderboard (Planilha1) was refreshed from the source data: every score
# in column C got its full floating-point precision back (previously rounded
# for display purposes), a new row was inserted at rank 16 ("Robson") and a
# new row was appended at the end as rank 36 ("Pedro André"). All the other
# rows keep the same id/name/score/rank/link, just shifted down by one where
# the new row was inserted above them.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Insert a new row at 17 ("Robson") - shifts old rows 17-35 down to 18-36,
# carrying the autofilter / used range along with it.
$ws.Rows("17:17").Insert(-4121)

# ---- Re-write the data rows (A id, B name, C score, D rank, E link) ----

$ws.Range("A2").Value = "6053d6d93bc68a306598cc8a"
$ws.Range("B2").Value = "Rogério Imai"
$ws.Range("C2").Value = 2779715780.1971998
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "MitoBR"

$ws.Range("A3").Value = "65e09a04a3eedf92402bc595"
$ws.Range("B3").Value = "Hugo Matos"
$ws.Range("C3").Value = 2364740766.4398003
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "BabyTux69"

$ws.Range("A4").Value = "61e484ca5aa1be001868f065"
$ws.Range("B4").Value = "HANTAROGAMER"
$ws.Range("C4").Value = 2354120786.6111999
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "HANTAROGAMER"

$ws.Range("A5").Value = "62117ac581c6a00035d1e7f5"
$ws.Range("B5").Value = "BillyTKD"
$ws.Range("C5").Value = 1647063840.4191999
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = "kzugpfwl"

$ws.Range("A6").Value = "623dfbf9584b270011716c7d"
$ws.Range("B6").Value = "ZeraTAL"
$ws.Range("C6").Value = 1202672016.6792002
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = "ZeraTAL"

$ws.Range("A7").Value = "61e852b4dc27dc001969efa3"
$ws.Range("B7").Value = "RKFox"
$ws.Range("C7").Value = 1153129978.0010002
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = "RKFox"

$ws.Range("A8").Value = "6095fca68a19d000196923bf"
$ws.Range("B8").Value = "Yarey"
$ws.Range("C8").Value = 968327492.69599998
$ws.Range("D8").Value = 7
$ws.Range("E8").Value = "YesoGengo"

$ws.Range("A9").Value = "655d47c4196526c7c27a44b0"
$ws.Range("B9").Value = "Wallace Jack"
$ws.Range("C9").Value = 648955510.75940001
$ws.Range("D9").Value = 8
$ws.Range("E9").Value = "Walljack"

$ws.Range("A10").Value = "62e828043bcdc2d82492325d"
$ws.Range("B10").Value = "Hugo Menezes"
$ws.Range("C10").Value = 590539773.57160008
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = "yokko"

$ws.Range("A11").Value = "5b456fbf4c0cfe0013fec57e"
$ws.Range("B11").Value = "Lion Fera"
$ws.Range("C11").Value = 579670695.75
$ws.Range("D11").Value = 10
$ws.Range("E11").Value = "jzlfzveq"

$ws.Range("A12").Value = "637c24d32bc15392f36cb7d3"
$ws.Range("B12").Value = "Fábio HK"
$ws.Range("C12").Value = 450714236.7712
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = "larjauxv"

$ws.Range("A13").Value = "65cdf66809b03f95ead3f9d2"
$ws.Range("B13").Value = "Rafael Manthy"
$ws.Range("C13").Value = 411764926.16729999
$ws.Range("D13").Value = 12
$ws.Range("E13").Value = "manthy"

$ws.Range("A14").Value = "62f6a9a38f0594d1cedf63bc"
$ws.Range("B14").Value = "Carlos Salomão"
$ws.Range("C14").Value = 370711107.0995
$ws.Range("D14").Value = 13
$ws.Range("E14").Value = "l6qv2mk6"

$ws.Range("A15").Value = "6644084c20155f520d9877cb"
$ws.Range("B15").Value = "Jhow"
$ws.Range("C15").Value = 289204784.35119998
$ws.Range("D15").Value = 14
$ws.Range("E15").Value = "JHOWBR21"

$ws.Range("A16").Value = "65d7a4bb18dbdd3be41d49af"
$ws.Range("B16").Value = "Ricardo"
$ws.Range("C16").Value = 252411594.55039999
$ws.Range("D16").Value = 15
$ws.Range("E16").Value = "ricardomann"

$ws.Range("A17").Value = "61ec9d04b499da001220f1f3"
$ws.Range("B17").Value = "Robson"
$ws.Range("C17").Value = 251229137.727
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = "kyqi8rx2"

$ws.Range("A18").Value = "61e6f8aedcdbc50011a301af"
$ws.Range("B18").Value = "Renan"
$ws.Range("C18").Value = 233891214.64499998
$ws.Range("D18").Value = 17
$ws.Range("E18").Value = "jukinha"

$ws.Range("A19").Value = "65de9e82a1e9f41193e2f6cc"
$ws.Range("B19").Value = "GUERDE"
$ws.Range("C19").Value = 205750999.59819999
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = "GUERDE"

$ws.Range("A20").Value = "6432faab10e5731e1f00ad18"
$ws.Range("B20").Value = "PUDIMZINHO"
$ws.Range("C20").Value = 138072167.338
$ws.Range("D20").Value = 19
$ws.Range("E20").Value = "lg9p91px"

$ws.Range("A21").Value = "614b9f90608bdd002791f31c"
$ws.Range("B21").Value = "Terráqueo"
$ws.Range("C21").Value = 114976184.516
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = "terraqueo"

$ws.Range("A22").Value = "6116f40168427b0029d5580f"
$ws.Range("B22").Value = "Julio Barboza"
$ws.Range("C22").Value = 108158853.3408
$ws.Range("D22").Value = 21
$ws.Range("E22").Value = "SidFillips"

$ws.Range("A23").Value = "64c7fd88d5e874c8f5a00a07"
$ws.Range("B23").Value = "Matheus75K"
$ws.Range("C23").Value = 103098592.90099999
$ws.Range("D23").Value = 22
$ws.Range("E23").Value = "lkr7fomr"

$ws.Range("A24").Value = "66281a5897a97ed50ab05355"
$ws.Range("B24").Value = "Marcus"
$ws.Range("C24").Value = 102258845.814
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = "thekrk420"

$ws.Range("A25").Value = "663af1422c820df814115447"
$ws.Range("B25").Value = "Zonorato"
$ws.Range("C25").Value = 66382476.8719
$ws.Range("D25").Value = 24
$ws.Range("E25").Value = "lvx9ene8"

$ws.Range("A26").Value = "663590bd491ae03c93bf4900"
$ws.Range("B26").Value = "Dyego"
$ws.Range("C26").Value = 56256615.1142
$ws.Range("D26").Value = 25
$ws.Range("E26").Value = "Dyegolimax"

$ws.Range("A27").Value = "60d5e41a8185d30018933209"
$ws.Range("B27").Value = "Pedro Ferreira"
$ws.Range("C27").Value = 53401230.1735
$ws.Range("D27").Value = 26
$ws.Range("E27").Value = "Pesaac"

$ws.Range("A28").Value = "602ac6c9221f0d0036212c71"
$ws.Range("B28").Value = "LeoDGyn"
$ws.Range("C28").Value = 52439319.64
$ws.Range("D28").Value = 27
$ws.Range("E28").Value = "kl6yb7ip"

$ws.Range("A29").Value = "625461969628d1006fc9ea79"
$ws.Range("B29").Value = "Gustavo"
$ws.Range("C29").Value = 51238979.5989
$ws.Range("D29").Value = 28
$ws.Range("E29").Value = "sccp_gu"

$ws.Range("A30").Value = "607ba063fef5b1001a278047"
$ws.Range("B30").Value = "Orias Régis"
$ws.Range("C30").Value = 45036484.105
$ws.Range("D30").Value = 29
$ws.Range("E30").Value = "Regis"

$ws.Range("A31").Value = "61849ec5bee39a0026d534e7"
$ws.Range("B31").Value = "Fellype"
$ws.Range("C31").Value = 33605904.3755
$ws.Range("D31").Value = 30
$ws.Range("E31").Value = "kvlslck9"

$ws.Range("A32").Value = "6487ab9f2b120de602eabfb1"
$ws.Range("B32").Value = "Gon Sotans"
$ws.Range("C32").Value = 31856096.1484
$ws.Range("D32").Value = 31
$ws.Range("E32").Value = "lithrrgz"

$ws.Range("A33").Value = "6021e1404cf33e00111c2f1c"
$ws.Range("B33").Value = "Axel"
$ws.Range("C33").Value = 29982580.704
$ws.Range("D33").Value = 32
$ws.Range("E33").Value = "kkxb6fy2"

$ws.Range("A34").Value = "61f825f623ed52005e866ed3"
$ws.Range("B34").Value = "Lil Taro"
$ws.Range("C34").Value = 20157133.4176
$ws.Range("D34").Value = 33
$ws.Range("E34").Value = "kz30bi28"

$ws.Range("A35").Value = "6117fe5efc41af003a67eb8f"
$ws.Range("B35").Value = "Rafael Braga"
$ws.Range("C35").Value = 7540281.0574
$ws.Range("D35").Value = 34
$ws.Range("E35").Value = "ksc25irv"

$ws.Range("A36").Value = "6689b328679fcd0d46f06b1b"
$ws.Range("B36").Value = "Jonathan Soares"
$ws.Range("C36").Value = 2841504.14
$ws.Range("D36").Value = 35
$ws.Range("E36").Value = "lyamecmn"

$ws.Range("A37").Value = "5fe161d0c8e049004bfde7a1"
$ws.Range("B37").Value = "Pedro André"
$ws.Range("C37").Value = 244857.01500000001
$ws.Range("D37").Value = 36
$ws.Range("E37").Value = "kizekp4s"

# ---- Fix up per-cell formatting (column A/B use an Arial-Unicode "s=2"/
#      "s=3" style for some rows, matching the refreshed source data; a few
#      rows have no special formatting at all). Copy the desired format from
#      a stable donor cell so the shared style table stays exactly as it was.
$noStyle = $ws.Range("B2")     # plain, unformatted cell
$style1  = $ws.Range("A4")     # bold-ish "s=1" font style
$style2  = $ws.Range("A2")     # "s=2" Arial Unicode MS style
$style3  = $ws.Range("A10")    # "s=3" Arial Unicode MS + number-format style

$style2.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B2").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A3").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B3").PasteSpecial(-4122)
$style1.Copy()
$ws.Range("A4").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B4").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A5").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B5").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A6").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B6").PasteSpecial(-4122)
$style1.Copy()
$ws.Range("A7").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B7").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A8").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B8").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A9").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("B9").PasteSpecial(-4122)
$style3.Copy()
$ws.Range("A10").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B10").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A11").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B11").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A12").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B12").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A13").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B13").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("A14").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B14").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A15").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B15").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A16").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B16").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("A17").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B17").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A18").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B18").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("A19").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B19").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("A20").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("B20").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A21").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B21").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A22").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B22").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A23").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("B23").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A24").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B24").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A25").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B25").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A26").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B26").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A27").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B27").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A28").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("B28").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A29").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B29").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A30").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B30").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A31").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B31").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A32").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B32").PasteSpecial(-4122)
$style3.Copy()
$ws.Range("A33").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B33").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A34").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B34").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A35").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B35").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("A36").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B36").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("A37").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("B37").PasteSpecial(-4122)

# Column C keeps its "s=2" number style throughout, except rows 14, 17 and 37
# which are plain (unformatted) numbers.
$style2.Copy()
$ws.Range("C2").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C3").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C4").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C5").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C6").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C7").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C8").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C9").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C10").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C11").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C12").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C13").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("C14").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C16").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("C17").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C18").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C19").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C20").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C21").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C22").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C23").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C24").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C25").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C26").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C27").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C28").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C29").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C30").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C31").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C32").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C33").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C34").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C35").PasteSpecial(-4122)
$style2.Copy()
$ws.Range("C36").PasteSpecial(-4122)
$noStyle.Copy()
$ws.Range("C37").PasteSpecial(-4122)

$excel.CutCopyMode = 0
